$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (Finland, company_name "2")
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = -0.009349999999999999
$ws.Range("E2").Value = -0.0292
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2726.5
$ws.Range("L2").Value = 0.3012540743605326
$ws.Range("M2").Value = 58.0484
$ws.Range("N2").Value = 0.001734968423003237
$ws.Range("O2").Value = 0.02129044562626077
$ws.Range("P2").Value = 46.3484
$ws.Range("Q2").Value = 0.001385275226478649
$ws.Range("R2").Value = 0.0169992297817715
$ws.Range("S2").Value = 11.7
$ws.Range("T2").Value = 0.201555942971727
$ws.Range("U2").Value = 63958.4
$ws.Range("V2").Value = 1.911608319709247
$ws.Range("W2").Value = 0.1034305747368546
$ws.Range("X2").Value = 0.1691861193663548
$ws.Range("Y2").Value = -0.06575554462950023
$ws.Range("Z2").Value = 0.03287551113146899
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.03950298720896973
$ws.Range("AC2").Value = -0.03950298720896973
$ws.Range("AD2").Value = 281067.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 281067.4
$ws.Range("AG2").Value = 217109
$ws.Range("AH2").Value = 0.8936241377084768
$ws.Range("AI2").Value = 0.8793746611230332
$ws.Range("AJ2").Value = 0.8664711899297154
$ws.Range("AK2").Value = 0.8491982625599666
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# ---------------------------------------------------------------------------
# Row 3 (Finland, Nordea Bank Abp)
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = -0.09080000000000001
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2692.6
$ws.Range("L3").Value = 0.3029477947794779
$ws.Range("M3").Value = 43.4
$ws.Range("N3").Value = 0.001313006937154163
$ws.Range("O3").Value = 0.01611825001856941
$ws.Range("P3").Value = 31.7
$ws.Range("Q3").Value = 0.0009590396292116814
$ws.Range("R3").Value = 0.01177300750204264
$ws.Range("S3").Value = 11.7
$ws.Range("T3").Value = 0.2695852534562212
$ws.Range("U3").Value = 63300.8
$ws.Range("V3").Value = 1.915078099709868
$ws.Range("W3").Value = 0.08097923339498654
$ws.Range("X3").Value = 0.1989051040777692
$ws.Range("Y3").Value = -0.1179258706827827
$ws.Range("Z3").Value = 0.03251005972004332
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.03703177611460236
$ws.Range("AC3").Value = -0.03703177611460236
$ws.Range("AD3").Value = 278885.6
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 278885.6
$ws.Range("AG3").Value = 215584.8
$ws.Range("AH3").Value = 0.894037465598297
$ws.Range("AI3").Value = 0.8794565188941548
$ws.Range("AJ3").Value = 0.8670605179322447
$ws.Range("AK3").Value = 0.8493927556188573
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# ---------------------------------------------------------------------------
# Row 4 (Finland, Alandsbanken Abp)
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = 0.0217
$ws.Range("E4").Value = 0.03240000000000001
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 33.9
$ws.Range("L4").Value = 0.2086153846153846
$ws.Range("M4").Value = 14.6484
$ws.Range("N4").Value = 0.03625841584158416
$ws.Range("O4").Value = 0.4321061946902655
$ws.Range("P4").Value = 14.6484
$ws.Range("Q4").Value = 0.03625841584158416
$ws.Range("R4").Value = 0.4321061946902655
$ws.Range("U4").Value = 657.6
$ws.Range("V4").Value = 1.627722772277228
$ws.Range("W4").Value = 0.1258819160787226
$ws.Range("X4").Value = 0.1394671346549404
$ws.Range("Y4").Value = -0.01358521857621778
$ws.Range("Z4").Value = 0.08535560458031305
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.0419741983033371
$ws.Range("AC4").Value = -0.0419741983033371
$ws.Range("AD4").Value = 2181.8
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 2181.8
$ws.Range("AG4").Value = 1524.2
$ws.Range("AH4").Value = 0.8437620852347436
$ws.Range("AI4").Value = 0.8690352903688361
$ws.Range("AJ4").Value = 0.7904781661653356
$ws.Range("AK4").Value = 0.8225580140313007
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
